$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 105, shifting existing rows 105:108 down to 106:109
$ws.Rows.Item(105).Insert()

# Populate new row 105 with values, mirroring the layout of the surrounding rows
$ws.Cells.Item(105, 1).Value = 1
$ws.Cells.Item(105, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(105, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(105, 4).Value = 45239
$ws.Cells.Item(105, 5).Value = 15
$ws.Cells.Item(105, 6).Value = 100112009
$ws.Cells.Item(105, 7).Value = "Acelga"
$ws.Cells.Item(105, 8).Value = "Sin especificar"
$ws.Cells.Item(105, 9).Value = "Segunda"
$ws.Cells.Item(105, 10).Value = 500
$ws.Cells.Item(105, 11).Value = 800
$ws.Cells.Item(105, 12).Value = 1000
$ws.Cells.Item(105, 13).Value = 900
$ws.Cells.Item(105, 14).Value = "$/atado 2,5 a 3 kilos"
$ws.Cells.Item(105, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(105, 16).Value = 300
$ws.Cells.Item(105, 17).Value = 3
$ws.Cells.Item(105, 18).Value = "Hortaliza"
